# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) with newly-calculated/simulated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 8
    3  = 11
    4  = 7
    5  = 4
    6  = 7
    7  = 3
    8  = 10
    9  = 5
    10 = 7
    11 = 4
    12 = 3
    13 = 5
    14 = 10
    15 = 12
    16 = 4
    17 = 10
    18 = 3
    19 = 4
    20 = 8
    21 = 5
    22 = 6
    23 = 8
    24 = 14
    25 = 7
    26 = 5
    27 = 6
    28 = 5
    29 = 6
    30 = 5
    31 = 6
    32 = 5
    33 = 4
    34 = 4
    35 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
